$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Status updates (Open -> Fixed) for existing rows ---
$ws.Range("D3").Value = "Fixed"
$ws.Range("D5").Value = "Fixed"

# --- New rows 11 and 12 (written first so new shared strings land in the
#     same index order as the target workbook: 21, 22, then 23) ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Assets Using Record"
$ws.Range("D11").Value = "Open"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "mainWindowViewModel到处都是，要不要弄成全局的？"

# --- Row 5 gains a Description cell (rich text, "知道" bolded) + taller row ---
$richText = "在viewmodel中添加一个Save As命令，对应UI上一个按钮。按下这个按钮，即执行这个命令：" + [char]10 + "1. 新建一个workspace，其model的值，是从原model深度复制过来的" + [char]10 + "2. 将这个workspace添加到mainWindowViewModel" + [char]10 + "这里引发一个问题(ID 11)是，workspace (Both All and single)需要“知道”mainWindowViewModel"
$ws.Range("C5").Value = $richText
$ws.Range("C5").Characters(183, 2).Font.Bold = $true
$ws.Rows.Item(5).RowHeight = 60

# --- Selection moves to F5 ---
$ws.Range("F5").Select()
